# Update handback datetime values (Correspond Handoff Datetime / Correspond
# Handback DateTime) for the zh-cn and de-de sheets, simulating a re-run of
# the handback status report with newer timestamps.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-18 20:20:43"
$wsZh.Range("E5").Value = "2016-03-18 20:20:43"
$wsZh.Range("H2").Value = "2016-03-18 20:21:10"
$wsZh.Range("H5").Value = "2016-03-18 20:21:10"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-18 20:20:48"
$wsDe.Range("E5").Value = "2016-03-18 20:20:48"
$wsDe.Range("H2").Value = "2016-03-18 20:21:20"
$wsDe.Range("H5").Value = "2016-03-18 20:21:20"
